$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"5.346288662278841e-20"

$ws.Range("B3").Value = 0.4250667495367802
$ws.Range("C3").Value = 0.3698035559332329

$ws.Range("B4").Value = 0.02283686003961818
$ws.Range("C4").Value = [double]"2.719181816842377e-18"

$ws.Range("B5").Value = [double]"2.059984127722458e-17"
$ws.Range("C5").Value = [double]"4.683478700650488e-17"

$ws.Range("B6").Value = 0.3260486671197617
$ws.Range("C6").Value = 0.3176571920514362

$ws.Range("B7").Value = [double]"1.965116437629977e-18"
$ws.Range("C7").Value = 0.01762138445562067

$ws.Range("B8").Value = [double]"2.727446090158847e-19"

$ws.Range("B9").Value = 0.03395846391654626
$ws.Range("C9").Value = 0.05108086275464575

$ws.Range("B10").Value = 0.03464765278391641
$ws.Range("C10").Value = 0.07663993609146461

$ws.Range("B11").Value = 0.09670723458172271
$ws.Range("C11").Value = [double]"3.993166304001589e-18"

$ws.Range("C12").Value = [double]"1.492601215055902e-17"

$ws.Range("B13").Value = 0.0607343720216545
$ws.Range("C13").Value = 0.1671970687136
